# Updates the cryptos list (price + 1h volume columns), matching the
# GitHub Actions bot refresh commit "Updated cryptos list on
# Wed Jul 12 19:53:49 UTC 2023 with GitHub Actions".
#
# Column D ("Price") holds values that look numeric (e.g. "1.001",
# "244.00") but must stay TEXT, exactly like the rest of the sheet
# (the thousands separators even produce strings like "30.332.43" that
# aren't valid numbers at all). Plainly assigning .Value to a numeric-
# looking string lets Excel auto-coerce it to a Number, so each Price
# write is bracketed with a text NumberFormat ("@") to force String
# storage, then the cell style is reset to "Normal" so the cell keeps
# its original (unstyled/General) look, matching the source sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "30.332.43"
$ws.Range("E2").Value = "  -0.78%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.874.86"
$ws.Range("E3").Value = "  +0.01%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "244.00"
$ws.Range("E5").Value = "  -1.50%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.19%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.4684"
$ws.Range("E7").Value = "  -1.10%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.2869"
$ws.Range("E8").Value = "  -1.02%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("D9") "0.06426"
$ws.Range("E9").Value = "  -0.77%  "

# Row 10 - Solana
Set-TextValue $ws.Range("D10") "21.95"
$ws.Range("E10").Value = "  -0.21%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.07787"
$ws.Range("E11").Value = "  +0.77%  "

# Row 12 - WrappedEther
Set-TextValue $ws.Range("D12") "1.883.36"
$ws.Range("E12").Value = "  +0.51%  "

# Row 13 - Litecoin
Set-TextValue $ws.Range("D13") "95.03"
$ws.Range("E13").Value = "  -1.03%  "

# Row 14 - Polygon
Set-TextValue $ws.Range("D14") "0.7207"
$ws.Range("E14").Value = "  -2.83%  "

# Row 15 - Polkadot
Set-TextValue $ws.Range("D15") "5.135"
$ws.Range("E15").Value = "  -0.59%  "

# Row 16 - BitcoinCash
Set-TextValue $ws.Range("D16") "278.61"
$ws.Range("E16").Value = "  +1.52%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "30.325.50"
$ws.Range("E17").Value = "  -1.04%  "

# Row 18 - Avalanche
Set-TextValue $ws.Range("D18") "12.95"
$ws.Range("E18").Value = "  -2.36%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.18%  "

# Row 20 & 21 swap places (ShibaInu now ranks above WrappedliquidstakedEther2.0)
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D20") "0.000007393"
$ws.Range("E20").Value = "  -1.00%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D21") "2.130.62"
$ws.Range("E21").Value = "  +0.71%  "

# Row 22 - BinanceUSD
Set-TextValue $ws.Range("D22") "1.001"
$ws.Range("E22").Value = "  +0.05%  "

# Row 23 - Uniswap
Set-TextValue $ws.Range("D23") "5.212"
$ws.Range("E23").Value = "  +0.26%  "

# Row 24 - Chainlink
Set-TextValue $ws.Range("D24") "6.239"
$ws.Range("E24").Value = "  +1.10%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "163.57"
$ws.Range("E25").Value = "  -0.82%  "

# Row 26 - Cosmos
Set-TextValue $ws.Range("D26") "9.012"
$ws.Range("E26").Value = "  -2.03%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "18.64"
$ws.Range("E27").Value = "  -0.14%  "

# Row 28 - LidoDAOToken
Set-TextValue $ws.Range("D28") "1.878"
$ws.Range("E28").Value = "  -1.31%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  -1.17%  "

# Row 30 - Stellar
Set-TextValue $ws.Range("D30") "0.09567"
$ws.Range("E30").Value = "  -3.73%  "

# Row 31 - PancakeSwap
Set-TextValue $ws.Range("D31") "1.466"
$ws.Range("E31").Value = "  -2.80%  "

# Row 32 - Filecoin
Set-TextValue $ws.Range("D32") "4.203"
$ws.Range("E32").Value = "  -0.84%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D33") "4.087"
$ws.Range("E33").Value = "  +0.06%  "

# Row 34 - Hedera
Set-TextValue $ws.Range("D34") "0.04805"
$ws.Range("E34").Value = "  +0.76%  "

# Row 35 - ARBITRUM
Set-TextValue $ws.Range("D35") "1.118"
$ws.Range("E35").Value = "  -0.09%  "

# Row 36 - ImmutableX
Set-TextValue $ws.Range("D36") "0.6849"
$ws.Range("E36").Value = "  -0.97%  "

# Row 37 - HuobiToken
Set-TextValue $ws.Range("D37") "2.709"
$ws.Range("E37").Value = "  -0.32%  "

# Row 38 - VeChain
Set-TextValue $ws.Range("D38") "0.01866"
$ws.Range("E38").Value = "  +1.06%  "

# Row 39 - MXToken
Set-TextValue $ws.Range("D39") "2.811"
$ws.Range("E39").Value = "  +2.15%  "

# Row 40 - FraxShare
Set-TextValue $ws.Range("D40") "6.226"
$ws.Range("E40").Value = "  -0.62%  "

# Row 41 - Aave
Set-TextValue $ws.Range("D41") "74.34"
$ws.Range("E41").Value = "  +1.67%  "

# Row 42 - RenderToken
Set-TextValue $ws.Range("D42") "1.931"
$ws.Range("E42").Value = "  -1.74%  "

# Row 43 - TheSandbox
Set-TextValue $ws.Range("D43") "0.4216"
$ws.Range("E43").Value = "  +1.40%  "

# Row 44 - PaxDollar (unchanged)

# Row 45 - TrustWalletToken
Set-TextValue $ws.Range("D45") "0.8226"
$ws.Range("E45").Value = "  -1.29%  "

# Row 46 - Quant
Set-TextValue $ws.Range("D46") "100.66"
$ws.Range("E46").Value = "  -0.63%  "

# Row 47 - EnergySwap
Set-TextValue $ws.Range("D47") "9.527"
$ws.Range("E47").Value = "  +1.86%  "

# Row 48 - Elrond
Set-TextValue $ws.Range("D48") "34.99"
$ws.Range("E48").Value = "  -0.98%  "

# Row 49 - Aptos
Set-TextValue $ws.Range("D49") "6.911"
$ws.Range("E49").Value = "  -0.87%  "

# Row 50 - Maker
Set-TextValue $ws.Range("D50") "894.45"
$ws.Range("E50").Value = "  -2.12%  "

# Row 51 - Cronos
Set-TextValue $ws.Range("D51") "0.05722"
$ws.Range("E51").Value = "  +1.00%  "
